# Update loading_percent results for the 380 kV case (Case_5_204).
# Columns in sheetData order: B, C, D, E, F, I, J (G and H are intentionally
# absent/empty in this sheet and must stay untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columnIndexes = @(2, 3, 4, 5, 6, 9, 10)  # B, C, D, E, F, I, J

$newValues = @"
2,19.41746662583446,16.72799125855228,6.596485066606053,12.54334689024362,45.95126704239596,29.16451583828837,10.34759226140396
3,18.95060772146591,16.24472651386525,6.599381963313015,12.52409565201726,45.68090976713003,29.12847201596372,10.35653013457574
4,18.6640933686027,15.94628212193728,6.601847285039993,12.51497466854452,45.52828976688078,29.11418353287765,10.36396920002335
5,18.54756221238129,15.82444127546142,6.603025774198867,12.51193858857182,45.46949389641838,29.11032478413101,10.36749034192591
6,18.52823162319703,15.80420259586303,6.603231996308753,12.51147561726947,45.45993697096499,29.10980238550904,10.3681045680599
7,18.66252061825084,15.94463954920682,6.601862473065524,12.51493096408192,45.52748302721311,29.11412355259144,10.3640147061251
8,19.25657748970488,16.56184055656844,6.597342129389307,12.53615017023912,45.85529265750304,29.1504558123525,10.35026852355528
9,20.4145735254694,17.74990608433525,6.593873206090767,12.59907010935646,46.60253890635854,29.2842187317176,10.33883335615761
10,21.25092261244811,18.59854820119033,6.594543222289103,12.65811267732858,47.21245118728908,29.42088862709724,10.33994146324614
11,21.62630319079044,18.97737736105607,6.595531666670129,12.68770771906743,47.50244644421399,29.49142252029743,10.3425177144995
12,21.7675708748482,19.11964579542082,6.596002961247621,12.69930318441827,47.61399790117354,29.51933286165234,10.34379153494255
13,21.73718780021456,19.08906068610673,6.59589716687793,12.6967887012159,47.589897177586,29.51326851754334,10.34350392750885
14,21.63794384215997,18.9891064333999,6.595568504078243,12.6886539299416,47.51158944575344,29.49369467417652,10.34261653409835
15,21.5770350746362,18.9277229870329,6.595379780572255,12.68372158049097,47.46384778618653,29.48186141005066,10.34211182637193
16,21.2262752895448,18.5736325733484,6.594492269745779,12.65623321557216,47.19374607072794,29.41644722949594,10.33981482364709
17,21.00968407552692,18.35445042443637,6.594121978002726,12.64006741122453,47.03121358194225,29.37845881784961,10.33893666528554
18,20.88463678269178,18.22771075348627,6.593973424025465,12.63102715083646,46.93891411123136,29.35739619323556,10.33862660786137
19,20.84222200686599,18.18468820046666,6.593934231644831,12.62801071029375,46.90786858371586,29.35039999359683,10.33855511434687
20,21.03279032380307,18.37785338911575,6.59415473762376,12.64176163357396,47.04839328150329,29.38242127288083,10.3390099587562
21,21.66711924577224,19.01849871185998,6.595662418526532,12.69103280890251,47.53454375060969,29.49941142358787,10.3428690872436
22,22.07649414427947,19.43022626763575,6.597212499407926,12.72549598711765,47.86235921231823,29.58286762590068,10.34712969354208
23,21.8585252419842,19.21116248563492,6.596333958065401,12.7068971793595,47.68649808026277,29.53768638124926,10.34469661102913
24,21.02234561615961,18.36727516791679,6.594139726520636,12.64099488513201,47.04062277527626,29.3806274253616,10.33897621601874
25,20.1031377174659,17.43203923071044,6.59424125364424,12.57978338426797,46.38947200627929,29.24130302879465,10.34759226140396
"@

$rows = $newValues.Trim() -split "`n"
foreach ($line in $rows) {
    $fields = $line.Trim() -split ","
    $rowNum = [int]$fields[0]
    for ($i = 0; $i -lt $columnIndexes.Length; $i++) {
        $colIdx = $columnIndexes[$i]
        $value = [double]$fields[$i + 1]
        $ws.Cells.Item($rowNum, $colIdx).Value = $value
    }
}

